$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- VampireAmbiance / HellhoundAmbiance / PlayerAmbiance marked Complete, notes cleared ---
# (clear the old "Having trouble..." notes first so the new notes below reuse/compact the table
#  in the same order the original author produced them)
$ws.Range("G2").Copy()
$ws.Range("G21:G23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G21").Value = "Complete"
$ws.Range("G22").Value = "Complete"
$ws.Range("G23").Value = "Complete"
$ws.Range("H21").ClearContents()
$ws.Range("H22").ClearContents()
$ws.Range("H23").ClearContents()

# --- New general notes added in rows 28 and 29 (row 27 left blank) ---
$ws.Range("H28").Value = "Going to figure out how to pause"

# --- TutorialDialog note updated (pitch change fix note) ---
$ws.Range("H13").Value = "Fixed the one channel issue. My fault because of a setting I did with ReaEQ. I'm still not sure how much I like the sound though. "

$ws.Range("H29").Value = "Maybe add enemy health paramater to enemy ambiance with spawn rates of scatterers"

# --- Player health parameter notes added to SpookyMusic / BossMusic rows ---
$ws.Range("H18").Value = "Player Health parameter is added. Will probably need some tweaking yet"
$ws.Range("H19").Value = "Player Health parameter is added. Will probably need some tweaking yet"

# --- Update selection cursor position ---
$ws.Range("H21").Select()
